$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.999.89"
$ws.Range("E2").Value = "  -2.26%  "
$ws.Range("D3").Value = "2.590.11"
$ws.Range("E3").Value = "  -4.53%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'554.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("D6").Value = "'155.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  +0.69%  "
$ws.Range("D9").Value = "'0.105"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.55%  "
$ws.Range("D10").Value = "'0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.35%  "
$ws.Range("D11").Value = "'5.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.47%  "
$ws.Range("D12").Value = "'0.365"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").Value = "3.046.91"
$ws.Range("E13").Value = "  -4.66%  "
$ws.Range("D14").Value = "'25.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.83%  "
$ws.Range("D15").Value = "61.919.12"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "'0.0000144"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.59%  "
$ws.Range("D17").Value = "2.589.19"
$ws.Range("E17").Value = "  -4.68%  "
$ws.Range("D18").Value = "'11.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.13%  "
$ws.Range("D19").Value = "'4.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.45%  "
$ws.Range("D20").Value = "'339.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.44%  "
$ws.Range("D21").Value = "'6.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.16%  "
$ws.Range("D22").Value = "'0.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "'0.499"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.41%  "
$ws.Range("D24").Value = "'62.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'8.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").Value = "0.0₃0839"
$ws.Range("E28").Value = "  -5.63%  "
$ws.Range("D29").Value = "'1.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.79%  "
$ws.Range("D30").Value = "'7.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.46%  "
$ws.Range("D31").Value = "'1.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.95%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "'160.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.40%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").Value = "'0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "'19.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.73%  "
$ws.Range("D35").Value = "'4.70"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.66%  "
$ws.Range("D36").Value = "'1.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.87%  "
$ws.Range("D37").Value = "'1.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").Value = "'339.69"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.86%  "
$ws.Range("D39").Value = "'6.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("D40").Value = "'0.897"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.94%  "
$ws.Range("D41").Value = "'3.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.16%  "
$ws.Range("D42").Value = "'37.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.42%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'20.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.87%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.135.57"
$ws.Range("E45").Value = "  +1.59%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.608"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'19.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.87%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "'10.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "'0.0548"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.18%  "
$ws.Range("D50").Value = "'0.0966"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("D51").Value = "'0.0241"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.84%  "
